$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 118 - this shifts the existing row 118 (and its data)
# down to row 119, matching the diff where the old 44544 record moves to
# row 119 and a brand-new 44595 record takes the old row 118's position.
$ws.Rows.Item(118).Insert()

# Populate the newly inserted row 118 with the new weekly record.
$ws.Range("A118").Value = 8
$ws.Range("B118").Value = "Terminal La Palmera de La Serena"
$ws.Range("C118").Value = "Coquimbo"
$ws.Range("D118").Value = 44595
$ws.Range("E118").Value = 4
$ws.Range("F118").Value = 100112040
$ws.Range("G118").Value = "Cilantro"
$ws.Range("H118").Value = "Sin especificar"
$ws.Range("I118").Value = "Primera"
$ws.Range("J118").Value = 2400
$ws.Range("K118").Value = 2500
$ws.Range("L118").Value = 2800
$ws.Range("M118").Value = 2650
$ws.Range("N118").Value = '$/atado 1 a 1,5 kilos'
$ws.Range("O118").Value = "Provincia del Elquí"
$ws.Range("P118").Value = 1767
$ws.Range("Q118").Value = 1.5
$ws.Range("R118").Value = "Hortaliza"
